# Using internal column widths in pptx writer tables (#9392)
#
# The table on slide 1 ("Content Placeholder 5") has two grid columns
# that were re-measured by the pptx writer, growing from 2501900 EMU
# (197 pt) to 2514600 EMU (198 pt) each. Resize the table's columns
# through the PowerPoint object model so the grid reflects the new,
# internally-computed widths.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the graphicFrame shape that hosts the 2-column / 1-row table.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
    }
}

$tbl = $tableShape.Table

# New column widths, in points (2514600 EMU / 12700 EMU-per-point = 198 pt).
$newColWidthPts = 2514600 / 12700

$tbl.Columns.Item(1).Width = $newColWidthPts
$tbl.Columns.Item(2).Width = $newColWidthPts
